$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("APOLLO")

# Row 7 - top summary row
$ws.Range("F7").Value = 424.45
$ws.Range("G7").Value = 437.5
$ws.Range("H7").Value = 421.15
$ws.Range("I7").Value = 430.25
$ws.Range("J7").Value = 431.3

# Row 9
$ws.Range("G9").Value = 428.25
$ws.Range("H9").Value = 419.25
$ws.Range("I9").Value = 422.5

# Row 10
$ws.Range("G10").Value = 430
$ws.Range("H10").Value = 422.05
$ws.Range("I10").Value = 429.65

# Row 11
$ws.Range("G11").Value = 433.25
$ws.Range("H11").Value = 429.5
$ws.Range("I11").Value = 430.4

# Row 12
$ws.Range("G12").Value = 433.25
$ws.Range("H12").Value = 428.6
$ws.Range("I12").Value = 433.15

# Row 13
$ws.Range("G13").Value = 433.5
$ws.Range("H13").Value = 429.55
$ws.Range("I13").Value = 430.05

# Row 14
$ws.Range("G14").Value = 431.85
$ws.Range("H14").Value = 429
$ws.Range("I14").Value = 431.5

# Row 15
$ws.Range("G15").Value = 433.5
$ws.Range("H15").Value = 431.45
$ws.Range("I15").Value = 433.3

# Row 16
$ws.Range("G16").Value = 437.15
$ws.Range("H16").Value = 433.1
$ws.Range("I16").Value = 436.2

# Row 17
$ws.Range("G17").Value = 437.5
$ws.Range("H17").Value = 432.4
$ws.Range("I17").Value = 433.85

# Row 18
$ws.Range("G18").Value = 435.2
$ws.Range("H18").Value = 430.85
$ws.Range("I18").Value = 431.9

# Row 19
$ws.Range("G19").Value = 431.9
$ws.Range("H19").Value = 429.1
$ws.Range("I19").Value = 430.6

# Row 20
$ws.Range("G20").Value = 431.95
$ws.Range("H20").Value = 430
$ws.Range("I20").Value = 430.5

# Row 21
$ws.Range("G21").Value = 431.7
$ws.Range("H21").Value = 429.45
$ws.Range("I21").Value = 431.05
